$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 89; everything from row 89 down shifts to row 90..127
$ws.Rows.Item(89).Insert()

# Populate the new row 89 with the new weekly record
$ws.Cells.Item(89, 1).Value = 11
$ws.Cells.Item(89, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(89, 3).Value = "Bíobío"
$ws.Cells.Item(89, 4).Value = 45097
$ws.Cells.Item(89, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(89, 5).Value = 8
$ws.Cells.Item(89, 6).Value = 100112012
$ws.Cells.Item(89, 7).Value = "Espinaca"
$ws.Cells.Item(89, 8).Value = "Sin especificar"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 130
$ws.Cells.Item(89, 11).Value = 7500
$ws.Cells.Item(89, 12).Value = 8000
$ws.Cells.Item(89, 13).Value = 7769
$ws.Cells.Item(89, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(89, 15).Value = "Región Metropolitana"
$ws.Cells.Item(89, 16).Value = 777
$ws.Cells.Item(89, 17).Value = 10
$ws.Cells.Item(89, 18).Value = "Hortaliza"
